# "starting to implement keywriter method to gui"
# Append new Keys rows (69-91) to the "Keys" worksheet, mirroring data that a
# (still half-built) key-writer feature fed into the sheet: mostly blank rows
# with two real records (rows 88 and 90) in between.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Keys")

function Set-TextCell {
    param(
        $Sheet,
        [int]$Row,
        [int]$Col,
        [string]$Text
    )
    $cell = $Sheet.Cells.Item($Row, $Col)
    # Leading apostrophe forces a literal/text entry (so purely-numeric
    # strings like "2344" stay text instead of becoming numbers, and a bare
    # apostrophe with nothing after it yields a real, non-blank empty text
    # cell instead of clearing it back out).
    $cell.Value = "'" + $Text
    # Drop back to the workbook's default style so no quote-prefix / number
    # formatting is left behind on the cell.
    $cell.Style = "Normal"
}

# Rows 69-87: blank placeholder rows (all four columns hold an empty string).
for ($r = 69; $r -le 87; $r++) {
    Set-TextCell $ws $r 1 ""
    Set-TextCell $ws $r 2 ""
    Set-TextCell $ws $r 3 ""
    Set-TextCell $ws $r 4 ""
}

# Row 88: a real key record.
Set-TextCell $ws 88 1 "2344"
Set-TextCell $ws 88 2 "tamb"
Set-TextCell $ws 88 3 "lewis"
Set-TextCell $ws 88 4 "n533"

# Row 89: blank placeholder row again.
Set-TextCell $ws 89 1 ""
Set-TextCell $ws 89 2 ""
Set-TextCell $ws 89 3 ""
Set-TextCell $ws 89 4 ""

# Row 90: another real key record.
Set-TextCell $ws 90 1 "4334"
Set-TextCell $ws 90 2 "lle"
Set-TextCell $ws 90 3 "lewis"
Set-TextCell $ws 90 4 "n544"

# Row 91: trailing blank placeholder row.
Set-TextCell $ws 91 1 ""
Set-TextCell $ws 91 2 ""
Set-TextCell $ws 91 3 ""
Set-TextCell $ws 91 4 ""
